# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" columns (F/G) for
# both locale sheets with hyperlinked file names, stamps the "Latest Handback
# DateTime" column (H), and flips the handoff status text now that the
# localized content has been handed back and is in sync with en-US.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# --- zh-cn sheet ---------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Range("H2").Value = "2016-03-17 06:55:26"
$zh.Range("H3").Value = "2016-03-17 06:55:26"

$zh.Hyperlinks.Add(
    $zh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/831ab2472ff2429e00f4ad220bd3518a4843392a/e2e/c4819184-d820-4285-9d56-2a88b94b71ec.md",
    "",
    "",
    "c4819184-d820-4285-9d56-2a88b94b71ec.md"
)
$zh.Hyperlinks.Add(
    $zh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad3736f796b0060db25d088f19f93abb7b3fc555/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/c4819184-d820-4285-9d56-2a88b94b71ec.01c4a83c3c3a3708942f231f54fdf1d717b02539.zh-cn.xlf",
    "",
    "",
    "c4819184-d820-4285-9d56-2a88b94b71ec.01c4a83c3c3a3708942f231f54fdf1d717b02539.zh-cn.xlf"
)
$zh.Hyperlinks.Add(
    $zh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/831ab2472ff2429e00f4ad220bd3518a4843392a/e2e/dbce4adc-4e83-424c-b58d-f6eb7074e473.md",
    "",
    "",
    "dbce4adc-4e83-424c-b58d-f6eb7074e473.md"
)
$zh.Hyperlinks.Add(
    $zh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad3736f796b0060db25d088f19f93abb7b3fc555/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/dbce4adc-4e83-424c-b58d-f6eb7074e473.e0bacc4b73d13fe6e6dd97867dcb2ae6f151f8ea.zh-cn.xlf",
    "",
    "",
    "dbce4adc-4e83-424c-b58d-f6eb7074e473.e0bacc4b73d13fe6e6dd97867dcb2ae6f151f8ea.zh-cn.xlf"
)

# --- de-de sheet ----------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Range("H2").Value = "2016-03-17 06:55:39"
$de.Range("H3").Value = "2016-03-17 06:55:39"

$de.Hyperlinks.Add(
    $de.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/831ab2472ff2429e00f4ad220bd3518a4843392a/e2e/c4819184-d820-4285-9d56-2a88b94b71ec.md",
    "",
    "",
    "c4819184-d820-4285-9d56-2a88b94b71ec.md"
)
$de.Hyperlinks.Add(
    $de.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fef93bace929ef1a784aa7e6d971a8984ab557b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/c4819184-d820-4285-9d56-2a88b94b71ec.01c4a83c3c3a3708942f231f54fdf1d717b02539.de-de.xlf",
    "",
    "",
    "c4819184-d820-4285-9d56-2a88b94b71ec.01c4a83c3c3a3708942f231f54fdf1d717b02539.de-de.xlf"
)
$de.Hyperlinks.Add(
    $de.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/831ab2472ff2429e00f4ad220bd3518a4843392a/e2e/dbce4adc-4e83-424c-b58d-f6eb7074e473.md",
    "",
    "",
    "dbce4adc-4e83-424c-b58d-f6eb7074e473.md"
)
$de.Hyperlinks.Add(
    $de.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fef93bace929ef1a784aa7e6d971a8984ab557b/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/dbce4adc-4e83-424c-b58d-f6eb7074e473.e0bacc4b73d13fe6e6dd97867dcb2ae6f151f8ea.de-de.xlf",
    "",
    "",
    "dbce4adc-4e83-424c-b58d-f6eb7074e473.e0bacc4b73d13fe6e6dd97867dcb2ae6f151f8ea.de-de.xlf"
)
